# Re-applies the refreshed crypto snapshot (price + 1h volume-change
# columns) captured by the Sat Sep 16 05:45:13 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new text. Price-looking values that would otherwise be auto-
# coerced to a Number by Excel (losing e.g. a trailing '.00') are
# prefixed with a literal apostrophe so they stay Text, same as the
# original inline-string cells.
$updates = @(
    @('D2', '26.670.03'),
    @('E2', '  -0.06%  '),
    @('D3', '1.646.38'),
    @('E3', '  +0.59%  '),
    @('D4', '''1.00'),
    @('E4', '  +0.21%  '),
    @('D5', '''216.03'),
    @('E5', '  +1.07%  '),
    @('E6', '  -0.87%  '),
    @('E7', '  +0.18%  '),
    @('E8', '  -0.39%  '),
    @('D9', '''0.0626'),
    @('E9', '  +0.32%  '),
    @('D10', '''19.35'),
    @('E10', '  +0.36%  '),
    @('D11', '''0.0844'),
    @('E11', '  -0.09%  '),
    @('D12', '1.876.56'),
    @('E12', '  +0.58%  '),
    @('D13', '''4.22'),
    @('E13', '  +2.88%  '),
    @('D14', '1.630.61'),
    @('E14', '  -0.42%  '),
    @('E15', '  +1.30%  '),
    @('D16', '''66.27'),
    @('E16', '  +4.22%  '),
    @('D17', '26.716.69'),
    @('E17', '  +0.08%  '),
    @('E18', '  +1.29%  '),
    @('D19', '''219.70'),
    @('E19', '  -0.21%  '),
    @('D20', '''1.00'),
    @('E20', '  +0.25%  '),
    @('E21', '  +1.61%  '),
    @('E22', '  +1.96%  '),
    @('E23', '  +1.07%  '),
    @('D24', '''2.12'),
    @('E24', '  +10.30%  '),
    @('D25', '''147.33'),
    @('E25', '  -0.59%  '),
    @('D26', '''1.00'),
    @('E26', '  +0.12%  '),
    @('E27', '  -0.61%  '),
    @('D28', '''7.12'),
    @('E28', '  +2.66%  '),
    @('D29', '''15.90'),
    @('E29', '  +2.42%  '),
    @('D30', '''0.0517'),
    @('E30', '  +0.72%  '),
    @('E31', '  +0.61%  '),
    @('E32', '  +2.30%  '),
    @('E33', '  +2.25%  '),
    @('D34', '1.287.84'),
    @('E34', '  +5.85%  '),
    @('E35', '  +1.58%  '),
    @('E36', '  +6.22%  '),
    @('E37', '  +0.45%  '),
    @('D38', '''0.528'),
    @('E38', '  +4.02%  '),
    @('E39', '  +1.71%  '),
    @('E40', '  +0.21%  '),
    @('E41', '  +1.83%  '),
    @('E42', '  -2.02%  '),
    @('E43', '  -0.12%  '),
    @('D44', '1.788.54'),
    @('E44', '  +0.74%  '),
    @('D45', '''93.77'),
    @('E45', '  +0.50%  '),
    @('D46', '''60.15'),
    @('E46', '  +9.43%  '),
    @('E47', '  +3.73%  '),
    @('D48', '''0.0517'),
    @('D49', '''7.85'),
    @('E49', '  +1.94%  '),
    @('D50', '''0.0978'),
    @('E50', '  +3.17%  '),
    @('E51', '  -0.61%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
